# Implements: "implementing some test cases like (products & locations)"
#
# 1. Update the "Sign up" test username from testuser26 -> testuser33
# 2. Resize the saved window view
# 3. Add two new worksheets at the end: "Locations" and "Products"

$wb = $excel.ActiveWorkbook

# --- 1. Update username on "Sign up" sheet (I2) -----------------------------
$signUp = $wb.Worksheets.Item("Sign up")
$signUp.Range("I2").Value = "testuser33"

# --- 2. Shrink the workbook window -------------------------------------------
$win = $excel.ActiveWindow
$win.Width = 19200
$win.Height = 6780

# --- 3. Add "Locations" sheet -------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$locations = $wb.Worksheets.Add($null, $lastSheet)
$locations.Name = "Locations"

$locations.Range("A1").Value = "location site link"
$locations.Range("B1").Value = "Location  site heading"
$locations.Range("A2").Value = "https://www.parasoft.com/solutions/"
$locations.Range("B2").Value = "Deliver High-Quality & Secure Software"

$locations.Columns.Item(1).ColumnWidth = 38.2857142857143
$locations.Columns.Item(2).ColumnWidth = 39.5714285714286

$locations.Hyperlinks.Add($locations.Range("A2"), "https://www.parasoft.com/solutions/", "", "https://www.parasoft.com/solutions/", "https://www.parasoft.com/solutions/")
$locations.Range("A2").Style = "Hyperlink"

$locations.Range("A1").Select()

# --- 4. Add "Products" sheet --------------------------------------------------
$products = $wb.Worksheets.Add($null, $locations)
$products.Name = "Products"

$products.Range("A1").Value = "Productss site link"
$products.Range("B1").Value = "Products site heading"
$products.Range("A2").Value = "https://www.parasoft.com/products/"
$products.Range("B2").Value = "Products"

$products.Columns.Item(1).ColumnWidth = 37.8571428571429
$products.Columns.Item(2).ColumnWidth = 22

$products.Hyperlinks.Add($products.Range("A2"), "https://www.parasoft.com/products/", "", "", "https://www.parasoft.com/products/")
$products.Range("A2").Style = "Hyperlink"

$signUp.Select()
